$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RuntimesChart")

# Day 5 (part 2) solved -- add new data point to the runtimes table.
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 0.0030532

# Update selection to match the newly-added row (mirrors Excel's default
# behavior of moving the active cell/selection to the row just edited).
$ws.Range("A7:B7").Select()
